$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.982.70'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.862.79'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9984'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9988'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5062'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.01%  '
$ws.Range('E8').Value = '  -0.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07142'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8832'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.63'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.877.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07563'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.299'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.18'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9986'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008415'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.30%  '
$ws.Range('E18').Value = '  -2.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9991'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.053.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.032'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.094.76'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.87%  '
$ws.Range('E24').Value = '  -1.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.842'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '147.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.098'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.72'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.671'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.702'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09044'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05135'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.033'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.153'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7291'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02039'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.035'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.460'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.076'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5279'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.533'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '115.83'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.283'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1468'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9985'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4611'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.942'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.565'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.53'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.91'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.61%  '
